$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Delete the useless "Content" (table-of-contents) slide.
#    It only contains a title "Content" plus three textboxes
#    "1. Our team" / "2.About the project" / "3.Used technologies" that
#    just duplicated the following slide titles, so it is removed outright.
# ---------------------------------------------------------------------------
$contentSlideIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Count -ge 1 -and $candidate.Shapes.Item(1).HasTextFrame) {
        if ($candidate.Shapes.Item(1).TextFrame.TextRange.Text -eq "Content") {
            $contentSlideIndex = $i
            break
        }
    }
}
if ($contentSlideIndex -eq -1) { $contentSlideIndex = 2 }
$p.Slides.Item($contentSlideIndex).Delete()

# ---------------------------------------------------------------------------
# 2. Clean up / reword the FAQ explanation on the "About the project" slide.
#    - merge the first two runs of the paragraph (no visible text change)
#    - replace the clumsy "thord ... oftenly asked" tail with a cleaner
#      "third option ... frequently asked questions and their answer."
# ---------------------------------------------------------------------------
$projectSlideIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Count -ge 1 -and $candidate.Shapes.Item(1).HasTextFrame) {
        if ($candidate.Shapes.Item(1).TextFrame.TextRange.Text -eq "About the project") {
            $projectSlideIndex = $i
            break
        }
    }
}
if ($projectSlideIndex -eq -1) { $projectSlideIndex = 3 }

$s = $p.Slides.Item($projectSlideIndex)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Merge the first two runs ("The project is ... menu with " + "3 options: ...
# your ") into a single run. Clear the first run, then re-write the now
# leading part of the (originally second) run so the combined text keeps
# the second run's formatting, exactly like Word/PowerPoint does when you
# backspace across a run boundary while typing.
$run1 = $tr.Characters(1, 93)
$run1.Text = ""
$run2 = $tr.Characters(1, 169)
$run2.Text = "The project is Employee of the month. When you run the project you will be given a menu with 3 options: Vote, Check prizes, FAQ. If you click the vote button you will be given a dropdown menu with all the workers for which you can vote for. When you choose your "

# Replace the tail (starting right after "votings") with the simplified text.
$tail = $tr.Characters(492, $tr.Length - 492 + 1)
$tail.Text = ". The third option (FAQ) is pretty simple. When you click it, you are given a list with some frequently asked questions and their answer."
